$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    ,("2023-10-08 10:42:29","hatespeech","def",10,327.9163126798707)
    ,("2023-10-08 10:42:29","hatespeech","setting2",10,383.8087181260709)
    ,("2023-10-08 10:42:29","hatespeech","setting1",10,377.6313738371288)
    ,("2023-10-08 10:42:29","hatespeech","setting5",10,240.49561539538)
    ,("2023-10-08 10:42:29","hatespeech","setting6",10,341.938791864331)
    ,("2023-10-08 10:42:29","hatespeech","setting3",10,639.6785897617817)
    ,("2023-10-08 10:42:30","hatespeech","setting4",10,353.3013859823798)
    ,("2023-10-08 10:42:30","hatespeech","def",20,649.5895881120642)
    ,("2023-10-08 10:42:30","hatespeech","setting7",10,343.3843535635364)
    ,("2023-10-08 10:42:30","hatespeech","setting6",20,645.5245582474438)
    ,("2023-10-08 10:42:30","hatespeech","setting1",20,610.2972736929402)
    ,("2023-10-08 10:42:30","hatespeech","setting7",20,642.8362613162063)
    ,("2023-10-08 10:42:30","hatespeech","setting3",20,832.7161550174501)
    ,("2023-10-08 10:42:31","hatespeech","setting4",20,620.5788724078254)
    ,("2023-10-08 10:42:31","hatespeech","setting8",20,569.9019976041708)
    ,("2023-10-08 10:42:31","hatespeech","setting5",20,577.5967446517869)
    ,("2023-10-08 10:42:31","hatespeech","setting8",10,318.353509847848)
    ,("2023-10-08 10:42:31","hatespeech","setting2",20,968.98071089194)
    ,("2023-10-08 10:42:31","hatespeech","setting8",30,843.6786016449133)
    ,("2023-10-08 10:42:32","hatespeech","setting4",40,1121.141390089806)
    ,("2023-10-08 10:42:32","hatespeech","def",40,1144.915338692704)
    ,("2023-10-08 10:42:32","hatespeech","setting5",30,828.1624513510438)
    ,("2023-10-08 10:42:32","hatespeech","setting2",30,970.921207006002)
    ,("2023-10-08 10:42:32","hatespeech","setting2",40,1192.21006259929)
    ,("2023-10-08 10:42:32","hatespeech","setting1",30,1003.866085537018)
    ,("2023-10-08 10:42:32","hatespeech","setting1",40,1203.24218690097)
    ,("2023-10-08 10:42:33","hatespeech","setting4",30,870.4034118661681)
    ,("2023-10-08 10:42:33","hatespeech","def",30,862.544538451858)
    ,("2023-10-08 10:42:33","hatespeech","setting7",30,858.2066333212081)
    ,("2023-10-08 10:42:33","hatespeech","setting6",30,858.0309278398564)
    ,("2023-10-08 10:42:33","hatespeech","setting3",30,1140.464697214593)
    ,("2023-10-08 10:42:34","hatespeech","setting3",40,1437.438961413374)
    ,("2023-10-08 10:42:34","hatespeech","setting5",40,1095.055386883123)
    ,("2023-10-08 10:42:34","hatespeech","setting6",40,1211.688976834227)
    ,("2023-10-08 10:42:34","hatespeech","setting7",40,1187.572505576021)
    ,("2023-10-08 10:42:34","hatespeech","setting8",40,1111.725865752564)
)

$startRow = 191
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Output "Done"
